$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Restructure rows: insert a new row at 7 so that the value that was in
# B6 (4.7300149999999999) moves down to B7, and a new "LogTidgi" row is
# created. This shifts old rows 7..27 to 8..28.
$ws.Rows.Item(7).Insert()

# --- Restructure rows: insert a new row at 10 (after the shift above, the
# old row 9 "GMTI" is now row 9) so a new empty "LogGMTI" row can sit
# between GMTI (row 9) and Pol (row 11, was row 9 pre-shift/row10 post one
# shift). This shifts rows 10..28 to 11..29.
$ws.Rows.Item(10).Insert()

# Clear the old B6 value (it moved to B7 conceptually; B6 must now be blank)
$ws.Range("B6").ClearContents()

# --- Column C: "updated name" duplicate/rename column ------------------
$ws.Range("C3").Value  = "Weiner"
$ws.Range("C4").Value  = "Mweiner"
$ws.Range("C5").Value  = "BalabanJ"
$ws.Range("C6").Value  = "Tigdi"
$ws.Range("C7").Value  = "LogTidgi"
$ws.Range("C8").Value  = "Xu"
$ws.Range("C9").Value  = "GMTI"
$ws.Range("C10").Value = "LogGMTI"
$ws.Range("C11").Value = "Pol"
$ws.Range("C12").Value = "DZ"
$ws.Range("C13").Value = "Ipc"
$ws.Range("C14").Value = "BertzCT"
$ws.Range("C15").Value = "Thara"
$ws.Range("C16").Value = "Tsch"
$ws.Range("C17").Value = "ZM1"
$ws.Range("C18").Value = "ZM2"
$ws.Range("C19").Value = "MZM1"
$ws.Range("C20").Value = "MZM2"
$ws.Range("C21").Value = "Qindex"
$ws.Range("C22").Value = "Platt"
$ws.Range("C23").Value = "diameter"
$ws.Range("C24").Value = "radius"
$ws.Range("C25").Value = "petitjean"
$ws.Range("C26").Value = "Sito"
$ws.Range("C27").Value = "Hato"
$ws.Range("C28").Value = "Geto"
$ws.Range("C29").Value = "Arto"

# --- Column F: "Update" reason column -----------------------------------
$ws.Range("F3").Value  = "Name update"
$ws.Range("F4").Value  = "Name update"
$ws.Range("F5").Value  = "Use RDKit source"
$ws.Range("F7").Value  = "Log of the Tidgi"
$ws.Range("F13").Value = "Use RDKIT source"
$ws.Range("F14").Value = "Use RDKit source"
$ws.Range("F23").Value = "Name update"
$ws.Range("F24").Value = "Name update"
$ws.Range("F25").Value = "Name update"

# --- New B7 value (moved down from former B6) ---------------------------
$ws.Range("B7").Value = 4.7300149999999999

# --- Selection / view ----------------------------------------------------
$ws.Range("J5").Select() | Out-Null

# --- Column widths ---------------------------------------------------------
# ColumnWidth (character units) -> stored xlsx width uses a pixel-grid
# rounding; use the inverse-computed values that reproduce the target
# stored widths as closely as possible.
$ws.Columns.Item(1).ColumnWidth = 13.166666666666666   # -> stored width 14
$ws.Columns.Item(5).ColumnWidth = 16.833333333333332   # -> stored width ~17.71
$ws.Columns.Item(6).ColumnWidth = 14.666666666666666   # -> stored width ~15.43

# --- Page setup : portrait orientation -----------------------------------
$ws.PageSetup.Orientation = 1
